{"js": "// \"Specified 'five' co-op terms\"\n// The EXPERIENCE section's first bullet reads \"Full-time co-ops, January 2013 \u2013\n// August 2015\"; the author clarified the count by turning it into\n// \"Five full-time co-ops, January 2013 \u2013 August 2015\". Word also re-anchors its\n// internal \"last edit\" (_GoBack) bookmark to the new edit location, so we move\n// it too (delete the old one, re-add it at the new insertion point).\n\nconst body = context.document.body;\n\n// Locate the exact phrase that needs to change. It occurs once, near the\n// start of the EXPERIENCE section (\"Full-time co-ops, ...\").\nconst matches = body.search(\"Full-time \", { matchCase: true, matchWholeWord: false });\nmatches.load(\"items\");\nawait context.sync();\n\nif (matches.items.length === 0) {\n  throw new Error(\"Could not find 'Full-time ' text to update.\");\n}\n\nconst target = matches.items[0];\ntarget.insertText(\"Five full-time \", \"Replace\");\nawait context.sync();\n\n// Word keeps a hidden \"_GoBack\" bookmark marking the location of the most\n// recent edit. Remove the old one (it sat after the \"advisory program\"\n// bullet from a prior edit) and drop a fresh one at the new edit spot, right\n// after the word \"Five\" that was just typed.\nconst doc = context.document;\nconst goBack = doc.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\nif (!goBack.isNullObject) {\n  doc.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nconst fiveMatches = body.search(\"Five\", { matchCase: true, matchWholeWord: false });\nfiveMatches.load(\"items\");\nawait context.sync();\n\nif (fiveMatches.items.length > 0) {\n  const afterFive = fiveMatches.items[0].getRange(\"End\");\n  afterFive.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"Specified 'five' co-op terms\"\n# The EXPERIENCE section's first bullet reads \"Full-time co-ops, January 2013 -\n# August 2015\"; the author clarified the count by turning it into\n# \"Five full-time co-ops, January 2013 - August 2015\". Word also re-anchors its\n# internal \"last edit\" (_GoBack) bookmark to the new edit location, so we move\n# it too (delete the old one, re-add it at the new insertion point).\n\n$d = $word.ActiveDocument\n\n# Locate the exact phrase that needs to change and replace it in place.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Full-time \"\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$found = $rng.Find.Execute()\nif ($rng.Find.Found) {\n  $rng.Text = \"Five full-time \"\n}\n\n# Word keeps a hidden \"_GoBack\" bookmark marking the location of the most\n# recent edit. Remove the old one (it sat after the \"advisory program\" bullet\n# from a prior edit) and drop a fresh one at the new edit spot, right after\n# the word \"Five\" that was just typed.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rng2 = $d.Content\n$rng2.Find.ClearFormatting()\n$rng2.Find.Text = \"Five\"\n$rng2.Find.MatchCase = $true\n$rng2.Find.MatchWholeWord = $false\n$null = $rng2.Find.Execute()\nif ($rng2.Find.Found) {\n  $rng2.Collapse(0)\n  $d.Bookmarks.Add(\"_GoBack\", $rng2)\n}\n"}
